$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 85.5428312652268
$ws.Range("D2").Value = 0.0000000000000000116359219675931

$ws.Range("B3").Value = 3267.31855575395

$ws.Range("B4").Value = 159.77401225967
$ws.Range("D4").Value = 0.0000000000000000000000000000360763509549953
